$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 53, pushing the existing rows 53-56 down to 54-57.
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with the new weekly record
# (same shape as the surrounding rows, new date / quality / volume / prices).
$ws.Range("A53").Value = 8
$ws.Range("B53").Value = "Terminal La Palmera de La Serena"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 44706
$ws.Range("E53").Value = 4
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100104
$ws.Range("H53").Value = "Frutos de pepita"
$ws.Range("I53").Value = 100104003
$ws.Range("J53").Value = "Membrillo"
$ws.Range("K53").Value = "Champion"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 16
$ws.Range("N53").Value = 280000
$ws.Range("O53").Value = 290000
$ws.Range("P53").Value = 285000
$ws.Range("Q53").Value = "$/bins (450 kilos)"
$ws.Range("R53").Value = "Región de O'Higgins"
$ws.Range("S53").Value = 633
$ws.Range("T53").Value = 450
